# Apply updated Betfair back/lay odds for Jogos_do_Dia workbook (2026-01-07)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.48
$ws.Range("F3").Value = 3.25
$ws.Range("G3").Value = 3.75
$ws.Range("H3").Value = 2.08
$ws.Range("I3").Value = 2.26
$ws.Range("K3").Value = 4.4
$ws.Range("L3").Value = 1.26
$ws.Range("N3").Value = 4.7
$ws.Range("P3").Value = 2.26
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 1.51
$ws.Range("S3").Value = 2.64
$ws.Range("T3").Value = 1.6
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.79
$ws.Range("AA3").Value = 32
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 11
$ws.Range("AF3").Value = 30
$ws.Range("AG3").Value = 17
$ws.Range("AH3").Value = 18.5
$ws.Range("AJ3").Value = 65
$ws.Range("AK3").Value = 40
$ws.Range("AL3").Value = 46
$ws.Range("AN3").Value = 29
$ws.Range("F4").Value = 4.4
$ws.Range("S4").Value = 4
$ws.Range("U4").Value = 1.84
$ws.Range("G5").Value = 980
$ws.Range("H5").Value = 2.12
$ws.Range("J5").Value = 1.03
$ws.Range("N5").Value = 1.29
$ws.Range("P5").Value = 1.29
$ws.Range("S5").Value = 1.05
$ws.Range("W5").Value = 1.41
$ws.Range("N6").Value = 1.1
$ws.Range("O6").Value = 1.2
$ws.Range("P6").Value = 1.24
$ws.Range("Q6").Value = 1.19
$ws.Range("S6").Value = 1.2
$ws.Range("J7").Value = 1.03
$ws.Range("N7").Value = 1.1
$ws.Range("S7").Value = 1.2
$ws.Range("H8").Value = 10.5
$ws.Range("I8").Value = 11
$ws.Range("J8").Value = 5.3
$ws.Range("K8").Value = 5.4
$ws.Range("T8").Value = 2.6
$ws.Range("V8").Value = 1.1
$ws.Range("W8").Value = 3.55
$ws.Range("Y8").Value = 27
$ws.Range("AA8").Value = 640
$ws.Range("AE8").Value = 290
$ws.Range("AO8").Value = 500
$ws.Range("F9").Value = 2.94
$ws.Range("G9").Value = 2.96
$ws.Range("H9").Value = 2.72
$ws.Range("I9").Value = 2.74
$ws.Range("N9").Value = 3.65
$ws.Range("O9").Value = 1.35
$ws.Range("P9").Value = 1.89
$ws.Range("R9").Value = 1.34
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.8
$ws.Range("U9").Value = 2.2
$ws.Range("V9").Value = 1.57
$ws.Range("Z9").Value = 17
$ws.Range("AK9").Value = 34
$ws.Range("AM9").Value = 100
$ws.Range("AN9").Value = 30
$ws.Range("H10").Value = 2.18
$ws.Range("I10").Value = 2.2
$ws.Range("Q10").Value = 1.84
$ws.Range("T10").Value = 1.71
$ws.Range("U10").Value = 2.36
$ws.Range("V10").Value = 1.83
$ws.Range("Z10").Value = 14
$ws.Range("AN10").Value = 32
$ws.Range("AO10").Value = 14.5
$ws.Range("H11").Value = 4.8
$ws.Range("I11").Value = 4.9
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 1.89
$ws.Range("Q11").Value = 2.1
$ws.Range("R11").Value = 1.33
$ws.Range("S11").Value = 3.85
$ws.Range("T11").Value = 1.95
$ws.Range("W11").Value = 2.1
$ws.Range("AG11").Value = 10.5
$ws.Range("AK11").Value = 20
$ws.Range("F12").Value = 3.45
$ws.Range("G12").Value = 3.5
$ws.Range("H12").Value = 2.34
$ws.Range("I12").Value = 2.36
$ws.Range("N12").Value = 3.95
$ws.Range("P12").Value = 1.98
$ws.Range("Q12").Value = 1.98
$ws.Range("U12").Value = 2.24
$ws.Range("V12").Value = 1.73
$ws.Range("W12").Value = 1.4
$ws.Range("Y12").Value = 10.5
$ws.Range("Z12").Value = 14.5
$ws.Range("AA12").Value = 30
$ws.Range("AB12").Value = 14
$ws.Range("AE12").Value = 24
$ws.Range("AF12").Value = 24
$ws.Range("AG12").Value = 14.5
$ws.Range("AI12").Value = 36
$ws.Range("AJ12").Value = 60
$ws.Range("AL12").Value = 48
$ws.Range("AN12").Value = 36
$ws.Range("AO12").Value = 18.5
$ws.Range("F13").Value = 2.22
$ws.Range("G13").Value = 2.26
$ws.Range("H13").Value = 3.45
$ws.Range("I13").Value = 3.5
$ws.Range("L13").Value = 1.35
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 4.5
$ws.Range("P13").Value = 2.16
$ws.Range("Q13").Value = 1.85
$ws.Range("R13").Value = 1.46
$ws.Range("S13").Value = 3.05
$ws.Range("T13").Value = 1.69
$ws.Range("U13").Value = 2.4
$ws.Range("V13").Value = 1.4
$ws.Range("W13").Value = 1.8
$ws.Range("X13").Value = 17.5
$ws.Range("AF13").Value = 14.5
$ws.Range("AM13").Value = 75
$ws.Range("AN13").Value = 14.5
$ws.Range("AO13").Value = 32
$ws.Range("F14").Value = 1.42
$ws.Range("G14").Value = 1.44
$ws.Range("H14").Value = 8
$ws.Range("I14").Value = 8.199999999999999
$ws.Range("J14").Value = 5.6
$ws.Range("K14").Value = 5.7
$ws.Range("M14").Value = 1.03
$ws.Range("P14").Value = 3.1
$ws.Range("R14").Value = 1.85
$ws.Range("S14").Value = 2.14
$ws.Range("T14").Value = 1.67
$ws.Range("W14").Value = 3.25
$ws.Range("X14").Value = 34
$ws.Range("AA14").Value = 220
$ws.Range("AC14").Value = 13
$ws.Range("AD14").Value = 29
$ws.Range("AE14").Value = 90
$ws.Range("AF14").Value = 11.5
$ws.Range("AG14").Value = 9.800000000000001
$ws.Range("AH14").Value = 19.5
$ws.Range("AI14").Value = 70
$ws.Range("AK14").Value = 13
$ws.Range("AM14").Value = 75
$ws.Range("AN14").Value = 4.6
$ws.Range("AO14").Value = 70
$ws.Range("F15").Value = 1.84
$ws.Range("G15").Value = 1.85
$ws.Range("N15").Value = 3.45
$ws.Range("P15").Value = 1.81
$ws.Range("U15").Value = 1.92
$ws.Range("W15").Value = 2.16
$ws.Range("Y15").Value = 16
$ws.Range("AE15").Value = 75
$ws.Range("AH15").Value = 22
$ws.Range("AM15").Value = 140
$ws.Range("H16").Value = 3.35
$ws.Range("I16").Value = 3.4
$ws.Range("N16").Value = 2.98
$ws.Range("O16").Value = 1.49
$ws.Range("Q16").Value = 2.48
$ws.Range("U16").Value = 1.9
$ws.Range("V16").Value = 1.41
$ws.Range("W16").Value = 1.61
$ws.Range("X17").Value = 9.6
$ws.Range("O18").Value = 1.23
$ws.Range("P18").Value = 2.36
$ws.Range("S18").Value = 2.74
$ws.Range("T18").Value = 2
$ws.Range("AJ18").Value = 340
$ws.Range("H19").Value = 1.73
$ws.Range("I19").Value = 1.74
$ws.Range("N19").Value = 4.4
$ws.Range("O19").Value = 1.27
$ws.Range("S19").Value = 3.05
$ws.Range("T19").Value = 1.81
$ws.Range("AC19").Value = 9
$ws.Range("F20").Value = 1.76
$ws.Range("G20").Value = 1.77
$ws.Range("H20").Value = 5.4
$ws.Range("I20").Value = 5.5
$ws.Range("K20").Value = 4.1
$ws.Range("N20").Value = 4.2
$ws.Range("V20").Value = 1.22
$ws.Range("W20").Value = 2.28
$ws.Range("Y20").Value = 19.5
$ws.Range("Z20").Value = 40
$ws.Range("AA20").Value = 130
$ws.Range("AC20").Value = 8.800000000000001
$ws.Range("AD20").Value = 19.5
$ws.Range("AE20").Value = 65
$ws.Range("AG20").Value = 9.6
$ws.Range("AJ20").Value = 17.5
$ws.Range("AK20").Value = 17
$ws.Range("AN20").Value = 10
